# The sheet's product/aliado list (PRODUCTOS!A2:D118) lost its last three
# rows (116-118): JOHANA CORREA/PAGO GIROS/..., MARIA ALEJANDRA/RETIROS BET
# PLAY/..., MARIA ALEJANDRA/OKI/... Clear the cell contents (keeping the
# existing cell formatting/styles) so the rows go blank, matching the
# shared-strings table shrinking from 202 to 197 unique entries once the
# now-unused strings ("PAGO GIROS", "RETIROS BET PLAY", "CORREDOR
# EMPRESARIAL", "OKI", "GRUPO RÉDITOS") drop out on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PRODUCTOS")

$ws.Range("A116:D118").ClearContents() | Out-Null

# Leave the selection where the author ended up after trimming the rows.
$ws.Range("D128").Select() | Out-Null
